$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.404.69"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.896.10"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.693"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0755"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0981"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.797"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.76%  "
$ws.Range("D15").Value = "2.172.46"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").Value = "1.870.80"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "35.453.30"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.50%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0603"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.26%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  +24.94%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -16.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.854"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0740"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0225"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.00%  "
$ws.Range("D45").Value = "1.327.34"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0810"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.16%  "
